# Automatic update of files.
# Update the "Förändrad" (Changed) date column C for rows 2-10
# from serial 45221 (2023-10-22) to serial 45224 (2023-10-25).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C10").Value = 45224
